$d = $word.ActiveDocument

# --- 1. Update the date/time in the document header ---
$d.Content.Find.Execute("June  21, 2021 (08:08:36 PM)", $true, $false, $false, $false, $false, $true, 1, $false, "June  21, 2021 (08:15:55 PM)", 2) | Out-Null

# --- 2. Append the new "Pushing Further (Optional)" section ---

function Insert-AtEnd($text) {
    $endPos = $d.Content.End
    $rng = $d.Range($endPos, $endPos)
    $rng.InsertAfter($text)
    $newEnd = $d.Content.End
    $len = $text.Length
    $startOfText = $newEnd - 1 - $len
    $endOfText = $newEnd - 1
    return @($startOfText, $endOfText)
}

$brk = Insert-AtEnd "`r"
$sectionStart = $brk[0]
$res = Insert-AtEnd "Pushing Further (Optional)"
$para1 = $d.Paragraphs.Last
$para1.Style = "Heading1"
$sectionEnd = $d.Content.End

$brk = Insert-AtEnd "`r"
$res = Insert-AtEnd "Start with two integer arrays with the following values:"
$para2 = $d.Paragraphs.Last
$para2.Style = "FirstParagraph"
$sectionEnd = $d.Content.End

$brk = Insert-AtEnd "`r"
$res = Insert-AtEnd "int"
$d.Range($res[0], $res[1]).Style = "DataTypeTok"
$res = Insert-AtEnd "[]"
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " left "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "="
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "{"
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "101"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "76"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "74"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "94"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "94"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "};"
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$null = Insert-AtEnd "`v"
$res = Insert-AtEnd "int"
$d.Range($res[0], $res[1]).Style = "DataTypeTok"
$res = Insert-AtEnd "[]"
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " right "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "="
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "{"
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "73"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "74"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "67"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "107"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "111"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "108"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "66"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "}"
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$para3 = $d.Paragraphs.Last
$para3.Style = "SourceCode"
$sectionEnd = $d.Content.End

$brk = Insert-AtEnd "`r"
$res = Insert-AtEnd "Implement statements to merge these two arrays, such that the resulting array contains the following values, in this order:"
$para4 = $d.Paragraphs.Last
$para4.Style = "FirstParagraph"
$sectionEnd = $d.Content.End

$brk = Insert-AtEnd "`r"
$res = Insert-AtEnd "101"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "76"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "74"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "94"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "94"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "73"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "74"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "67"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "107"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "111"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "108"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$res = Insert-AtEnd ","
$d.Range($res[0], $res[1]).Style = "OperatorTok"
$res = Insert-AtEnd " "
$d.Range($res[0], $res[1]).Style = "NormalTok"
$res = Insert-AtEnd "66"
$d.Range($res[0], $res[1]).Style = "DecValTok"
$para5 = $d.Paragraphs.Last
$para5.Style = "SourceCode"
$sectionEnd = $d.Content.End

$brk = Insert-AtEnd "`r"
$res = Insert-AtEnd "Do not use built-in array methods."
$para6 = $d.Paragraphs.Last
$para6.Style = "FirstParagraph"
$sectionEnd = $d.Content.End

# --- 3. Wrap the new section in a bookmark ---
$bookmarkRange = $d.Range($sectionStart, $sectionEnd)
$d.Bookmarks.Add("pushing-further-optional", $bookmarkRange) | Out-Null

Write-Output "Edit applied successfully"